$d = $word.ActiveDocument

$replacements = @(
    @{old = "16×59=944";  new = "35×46=1610"},
    @{old = "79×56=4424"; new = "92×69=6348"},
    @{old = "68×60=4080"; new = "28×15=420"},
    @{old = "71×35=2485"; new = "19×69=1311"},
    @{old = "17×66=1122"; new = "51×76=3876"},
    @{old = "42×88=3696"; new = "47×88=4136"},
    @{old = "60×50=3000"; new = "39×40=1560"},
    @{old = "98×67=6566"; new = "17×87=1479"},
    @{old = "79×61=4819"; new = "39×46=1794"},
    @{old = "68×77=5236"; new = "58×86=4988"},
    @{old = "97×41=3977"; new = "79×15=1185"},
    @{old = "70×95=6650"; new = "94×97=9118"},
    @{old = "82×75=6150"; new = "42×21=882"},
    @{old = "23×35=805";  new = "72×28=2016"},
    @{old = "35×77=2695"; new = "51×69=3519"},
    @{old = "94×25=2350"; new = "55×18=990"},
    @{old = "44×30=1320"; new = "52×44=2288"},
    @{old = "79×18=1422"; new = "35×75=2625"},
    @{old = "81×88=7128"; new = "41×91=3731"},
    @{old = "32×92=2944"; new = "57×73=4161"},
    @{old = "54×12=648";  new = "48×67=3216"},
    @{old = "53×66=3498"; new = "38×60=2280"},
    @{old = "31×29=899";  new = "95×83=7885"},
    @{old = "74×89=6586"; new = "80×24=1920"},
    @{old = "73×92=6716"; new = "77×72=5544"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
